$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / Row 12, columns A:C currently hold (服务编号 / 2372) and (服务门店 / 店名) respectively.
# This edit swaps which row holds which label+value pair (服务编号 moves to row 12,
# 服务门店 moves to row 11), carrying each label's formatting along with it.

$a11 = $ws.Range("A11")
$b11 = $ws.Range("B11")
$c11 = $ws.Range("C11")
$a12 = $ws.Range("A12")
$b12 = $ws.Range("B12")
$c12 = $ws.Range("C12")

# Capture original values before overwriting anything.
$a11Val = $a11.Value
$b11Val = $b11.Value
$a12Val = $a12.Value
$b12Val = $b12.Value

# --- values ---
$a11.Value = $a12Val
$b11.Value = $b12Val
$a12.Value = $a11Val
$b12.Value = $b11Val

# --- formatting: row 11 (A:C) now takes on the "服务门店" look, row 12 takes the "服务编号" look ---
# 服务门店 row formatting (no border flag; general number format)
$a11.Font.Name = "宋体"
$a11.Font.Size = 12
$a11.HorizontalAlignment = -4131
$a11.VerticalAlignment = -4108
$a11.Borders.LineStyle = -4142

$b11.Font.Name = "宋体"
$b11.Font.Size = 12
$b11.HorizontalAlignment = -4131
$b11.VerticalAlignment = -4108
$b11.NumberFormat = "General"
$b11.Borders.LineStyle = -4142

$c11.Font.Name = "宋体"
$c11.Font.Size = 12
$c11.HorizontalAlignment = -4131
$c11.VerticalAlignment = -4108
$c11.NumberFormat = "General"
$c11.Borders.LineStyle = -4142

# 服务编号 row formatting (border-apply flag set historically; general number format)
$a12.Font.Name = "宋体"
$a12.Font.Size = 12
$a12.HorizontalAlignment = -4131
$a12.VerticalAlignment = -4108
$a12.Borders.LineStyle = -4142

$b12.Font.Name = "宋体"
$b12.Font.Size = 12
$b12.HorizontalAlignment = -4131
$b12.VerticalAlignment = -4108
$b12.NumberFormat = "General"
$b12.Borders.LineStyle = -4142

$c12.Font.Name = "宋体"
$c12.Font.Size = 12
$c12.HorizontalAlignment = -4131
$c12.VerticalAlignment = -4108
$c12.NumberFormat = "General"
$c12.Borders.LineStyle = -4142

# --- defined names follow their labels to the new rows ---
$svcNo = $wb.Names.Item("服务编号")
$svcNo.RefersTo = "=浙江杭州三墩地铁站店!`$B`$12"

$svcStore = $wb.Names.Item("服务门店")
$svcStore.RefersTo = "=浙江杭州三墩地铁站店!`$B`$11"

# --- re-assert the merges for B11:C11 / B12:C12 in that order ---
$ws.Range("B11:C11").UnMerge() | Out-Null
$ws.Range("B12:C12").UnMerge() | Out-Null
$ws.Range("B11:C11").Merge() | Out-Null
$ws.Range("B12:C12").Merge() | Out-Null

Write-Host "swap complete"
